$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 107
$ws.Range("F9").Value = 866
$ws.Range("F11").Value = 1263
$ws.Range("F12").Value = 1508
$ws.Range("F15").Value = 337
$ws.Range("F16").Value = 1679
$ws.Range("F18").Value = 1088
$ws.Range("F22").Value = 1683
$ws.Range("F26").Value = 1187
$ws.Range("F27").Value = 307354
$ws.Range("F29").Value = 73
$ws.Range("F35").Value = 1127
$ws.Range("F36").Value = 1058
$ws.Range("F37").Value = 253
$ws.Range("F39").Value = 872
$ws.Range("F43").Value = 1094
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 264
$ws.Range("F6").Value = 4637
$ws.Range("F10").Value = 725
$ws.Range("F11").Value = 462
$ws.Range("F13").Value = 1058
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 264
$ws.Range("F4").Value = 4637
$ws.Range("F5").Value = 725
$ws.Range("F9").Value = 1058
$ws.Range("F11").Value = 866
$ws.Range("F15").Value = 1263
$ws.Range("F16").Value = 1508
$ws.Range("F20").Value = 337
$ws.Range("F22").Value = 1679
$ws.Range("F24").Value = 1088
$ws.Range("F29").Value = 1683
$ws.Range("F34").Value = 1187
$ws.Range("F37").Value = 73
$ws.Range("F41").Value = 1127
$ws.Range("F43").Value = 253
$ws.Range("F44").Value = 872
